# Updated cryptos list on Thu Dec  7 19:46:13 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for each coin row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily force the Price data cells (column D, rows 2-51) to Text format
# so numeric-looking price strings (e.g. "0.649", "231.14") are stored as text,
# matching the source data, and are not reinterpreted by Excel as numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.247.02"
$ws.Range("E2").Value = "  -2.00%  "
$ws.Range("D3").Value = "2.333.66"
$ws.Range("E3").Value = "  +2.92%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "0.649"
$ws.Range("E5").Value = "  +1.69%  "
$ws.Range("D6").Value = "231.14"
$ws.Range("E6").Value = "  +0.06%  "
$ws.Range("D7").Value = "65.21"
$ws.Range("E7").Value = "  +1.82%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.451"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").Value = "0.0950"
$ws.Range("E10").Value = "  -4.96%  "
$ws.Range("D11").Value = "56.70"
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("D12").Value = "26.58"
$ws.Range("E12").Value = "  -2.68%  "
$ws.Range("D13").Value = "2.680.87"
$ws.Range("E13").Value = "  +2.89%  "
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("D15").Value = "15.27"
$ws.Range("E15").Value = "  -3.31%  "
$ws.Range("E16").Value = "  +1.66%  "
$ws.Range("D17").Value = "0.838"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("D18").Value = "2.336.05"
$ws.Range("E18").Value = "  +2.79%  "
$ws.Range("D19").Value = "43.205.76"
$ws.Range("E19").Value = "  -1.57%  "
$ws.Range("D20").Value = "0.0₃0971"
$ws.Range("E20").Value = "  -3.65%  "
$ws.Range("D21").Value = "73.65"
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("D22").Value = "6.16"
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("D23").Value = "247.09"
$ws.Range("E23").Value = "  -2.32%  "
$ws.Range("E24").Value = "  +16.91%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("D26").Value = "2.41"
$ws.Range("E26").Value = "  -1.84%  "
$ws.Range("E27").Value = "  +0.86%  "
$ws.Range("D28").Value = "9.85"
$ws.Range("E28").Value = "  -2.67%  "
$ws.Range("D29").Value = "174.58"
$ws.Range("E29").Value = "  +1.85%  "
$ws.Range("D30").Value = "22.12"
$ws.Range("E30").Value = "  +5.71%  "
$ws.Range("D31").Value = "1.47"
$ws.Range("E31").Value = "  +4.94%  "
$ws.Range("E32").Value = "  -7.75%  "
$ws.Range("E33").Value = "  +0.04%  "
$ws.Range("D34").Value = "5.01"
$ws.Range("E34").Value = "  +4.47%  "
$ws.Range("D35").Value = "0.0681"
$ws.Range("E35").Value = "  -3.78%  "
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("D37").Value = "2.47"
$ws.Range("E37").Value = "  +6.20%  "
$ws.Range("E38").Value = "  -1.31%  "
$ws.Range("D39").Value = "3.57"
$ws.Range("E39").Value = "  -6.33%  "
$ws.Range("D40").Value = "0.0249"
$ws.Range("E40").Value = "  -3.76%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("E42").Value = "  +7.90%  "
$ws.Range("D43").Value = "17.82"
$ws.Range("E43").Value = "  +1.87%  "
$ws.Range("E44").Value = "  +6.01%  "
$ws.Range("D45").Value = "98.19"
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("E46").Value = "  -1.41%  "
$ws.Range("D47").Value = "4.38"
$ws.Range("E47").Value = "  +0.06%  "
$ws.Range("D48").Value = "0.0941"
$ws.Range("D49").Value = "1.432.03"
$ws.Range("E49").Value = "  -1.03%  "
$ws.Range("D50").Value = "9.83"
$ws.Range("E50").Value = "  -5.81%  "
$ws.Range("D51").Value = "0.000202"
$ws.Range("E51").Value = "  -10.50%  "

# Restore the original (default) formatting/style for the Price data cells now
# that the text values are safely stored as strings.
$ws.Range("D2:D51").ClearFormats()
